$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.400.49"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.646.56"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'598.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "'154.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.547"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "2.646.18"
$ws.Range("E10").Value = "  +7.91%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "'5.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").Value = "'28.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.55%  "
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "3.125.61"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "68.265.33"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "2.659.92"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").Value = "'365.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "'7.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").Value = "'4.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.04%  "
$ws.Range("E23").Value = "  +2.27%  "
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("D25").Value = "'74.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'9.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("D29").Value = "2.774.61"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "'575.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'8.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.79%  "
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("E35").Value = "  +4.96%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +5.94%  "
$ws.Range("D38").Value = "'159.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").Value = "'19.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("E43").Value = "  +10.34%  "
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("E45").Value = "  +3.69%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "'157.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "'3.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "'22.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.74%  "
